$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8228499293327332
$ws.Range("B1").Value = 2.10289454460144
$ws.Range("D1").Value = 1.357670426368713
$ws.Range("E1").Value = 0.4972657263278961
